$wb = $excel.ActiveWorkbook

# Sheet "FINE" (sheet1.xml) - add row 3
$wsFine = $wb.Worksheets.Item("FINE")
$wsFine.Range("A3").Value = 2023
$wsFine.Range("B3").Value = 4
$wsFine.Range("C3").Value = 0.26
$wsFine.Range("D3").Value = 0.38
$wsFine.Range("E3").Value = 0.49
$wsFine.Range("F3").Value = 0.6

# Sheet "COARSE" (sheet2.xml) - add row 3
$wsCoarse = $wb.Worksheets.Item("COARSE")
$wsCoarse.Range("A3").Value = 2023
$wsCoarse.Range("B3").Value = 4
$wsCoarse.Range("C3").Value = 0.34
$wsCoarse.Range("D3").Value = 0.52
